$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("B2", 0.9662035273968854)
    ,@("C2", 0.2922930264897445)
    ,@("E2", 0.2756424402917332)
    ,@("F2", 1.537737671043047)
    ,@("G2", 0.3935621180173499)
    ,@("H2", 0.5532028149622406)
    ,@("I2", 0.5576134873791112)
    ,@("J2", 0.0236025723335942)
    ,@("L2", 0.5634547036626856)
    ,@("N2", 1.104131231087905)
    ,@("O2", 1.837966045474914)
    ,@("B3", 0.874220390729306)
    ,@("C3", 0.2865049226391676)
    ,@("E3", 0.2761449441553516)
    ,@("F3", 1.53497284433999)
    ,@("G3", 0.3930109363187242)
    ,@("H3", 0.5564545984607818)
    ,@("I3", 0.5634891040949377)
    ,@("J3", 0.02211049803630161)
    ,@("L3", 0.5509346629821295)
    ,@("N3", 1.105731281217004)
    ,@("O3", 1.843355514392769)
    ,@("B4", 0.8177401074844113)
    ,@("C4", 0.2829418039676597)
    ,@("E4", 0.2765545275132091)
    ,@("F4", 1.534106556762012)
    ,@("G4", 0.3929754575098556)
    ,@("H4", 0.558711521131336)
    ,@("I4", 0.5674181894446235)
    ,@("J4", 0.0211887848383725)
    ,@("L4", 0.5434519898627173)
    ,@("N4", 1.107110924337988)
    ,@("O4", 1.847843016475082)
    ,@("B5", 0.7947252112458614)
    ,@("C5", 0.2814876351450408)
    ,@("E5", 0.276746895738313)
    ,@("F5", 1.533962765968404)
    ,@("G5", 0.393037097967337)
    ,@("H5", 0.5596967328744924)
    ,@("I5", 0.5691001149242005)
    ,@("J5", 0.0208118000793327)
    ,@("L5", 0.5404545185319876)
    ,@("N5", 1.107773281824898)
    ,@("O5", 1.849967919950615)
    ,@("B6", 0.7909037368324618)
    ,@("C6", 0.2812460447660499)
    ,@("E6", 0.2767803775473539)
    ,@("F6", 1.533951531919492)
    ,@("G6", 0.3930519271034711)
    ,@("H6", 0.5598642835745835)
    ,@("I6", 0.5693842759195107)
    ,@("J6", 0.02074911940643531)
    ,@("L6", 0.5399599271156319)
    ,@("N6", 1.107889321850372)
    ,@("O6", 1.850338646127341)
    ,@("B7", 0.8174297129217507)
    ,@("C7", 0.2829222011297929)
    ,@("E7", 0.2765570187046684)
    ,@("F7", 1.534103770146658)
    ,@("G7", 0.3929759808025182)
    ,@("H7", 0.5587245427925467)
    ,@("I7", 0.5674405454062708)
    ,@("J7", 0.02118370624467047)
    ,@("L7", 0.543411354854527)
    ,@("N7", 1.107119451303618)
    ,@("O7", 1.847870474450943)
    ,@("B8", 0.9344892911132092)
    ,@("C8", 0.2902992822738639)
    ,@("E8", 0.2757947623745842)
    ,@("F8", 1.536611922726678)
    ,@("G8", 0.3933091436923633)
    ,@("H8", 0.5542700247881029)
    ,@("I8", 0.5595726654029782)
    ,@("J8", 0.0230892749220537)
    ,@("L8", 0.5590954926920944)
    ,@("N8", 1.104600641858553)
    ,@("O8", 1.839579707750104)
    ,@("B9", 1.163954720148581)
    ,@("C9", 0.3046870694107895)
    ,@("E9", 0.2750996327390069)
    ,@("F9", 1.548120771461555)
    ,@("G9", 0.39637057631316)
    ,@("H9", 0.5475987103921653)
    ,@("I9", 0.5466958303009548)
    ,@("J9", 0.02678104428073169)
    ,@("L9", 0.5914644613876021)
    ,@("N9", 1.102802386415831)
    ,@("O9", 1.832677526453779)
    ,@("B10", 1.332411440460533)
    ,@("C10", 0.3152032411204004)
    ,@("E10", 0.2750738585163539)
    ,@("F10", 1.560589489332088)
    ,@("G10", 0.4000948644259381)
    ,@("H10", 0.5439538801095409)
    ,@("I10", 0.5387933033924526)
    ,@("J10", 0.0294650600215931)
    ,@("L10", 0.6162161637288506)
    ,@("N10", 1.103383347942412)
    ,@("O10", 1.83332189610141)
    ,@("B11", 1.409002693521302)
    ,@("C11", 0.3199741469382502)
    ,@("E11", 0.2751669087884032)
    ,@("F11", 1.567132501282799)
    ,@("G11", 0.4021109891213825)
    ,@("H11", 0.5425682821859965)
    ,@("I11", 0.5355370615763846)
    ,@("J11", 0.03067976571141884)
    ,@("L11", 0.6276844945414837)
    ,@("N11", 1.104058131730483)
    ,@("O11", 1.834858779922712)
    ,@("B12", 1.437998308893668)
    ,@("C12", 0.3217787690831813)
    ,@("E12", 0.2752171632282625)
    ,@("F12", 1.569735276056321)
    ,@("G12", 0.4029208396431585)
    ,@("H12", 0.542082740443675)
    ,@("I12", 0.5343527503583942)
    ,@("J12", 0.03113882174204718)
    ,@("L12", 0.6320569810889083)
    ,@("N12", 1.104372461315279)
    ,@("O12", 1.835619771225311)
    ,@("B13", 1.431753964341908)
    ,@("C13", 0.3213902035317773)
    ,@("E13", 0.2752056728595313)
    ,@("F13", 1.569169161632999)
    ,@("G13", 0.402744359231292)
    ,@("H13", 0.542185569544344)
    ,@("I13", 0.5346056435683266)
    ,@("J13", 0.0310399974328206)
    ,@("L13", 0.6311139731109279)
    ,@("N13", 1.104302153060473)
    ,@("O13", 1.835447914306997)
    ,@("B14", 1.411388348605726)
    ,@("C14", 0.3201226555770091)
    ,@("E14", 0.2751707425380658)
    ,@("F14", 1.567344127304168)
    ,@("G14", 0.4021766857897973)
    ,@("H14", 0.542527551807936)
    ,@("I14", 0.5354386500860855)
    ,@("J14", 0.03071755124565811)
    ,@("L14", 0.628043628707843)
    ,@("N14", 1.104082814581545)
    ,@("O14", 1.834917798690327)
    ,@("B15", 1.398912741988624)
    ,@("C15", 0.3193459787167683)
    ,@("E15", 0.2751513011054492)
    ,@("F15", 1.566242525149605)
    ,@("G15", 0.4018350132824224)
    ,@("H15", 0.5427421241165149)
    ,@("I15", 0.5359552419426095)
    ,@("J15", 0.03051992220822086)
    ,@("L15", 0.6261668094484207)
    ,@("N15", 1.103956114763946)
    ,@("O15", 1.834616404110363)
    ,@("B16", 1.327405034994229)
    ,@("C16", 0.3148911773984935)
    ,@("E16", 0.2750698816840895)
    ,@("F16", 1.560179399774867)
    ,@("G16", 0.3999695920974773)
    ,@("H16", 0.5440499125732572)
    ,@("I16", 0.5390129261034247)
    ,@("J16", 0.0293855478688414)
    ,@("L16", 0.6154708521873147)
    ,@("N16", 1.103347490153084)
    ,@("O16", 1.833246498104415)
    ,@("B17", 1.283525563077205)
    ,@("C17", 0.3121548738650262)
    ,@("E17", 0.2750467346213554)
    ,@("F17", 1.556682821202585)
    ,@("G17", 0.3989077380353194)
    ,@("H17", 0.5449219630132518)
    ,@("I17", 0.5409754962217725)
    ,@("J17", 0.02868802280647742)
    ,@("L17", 0.6089624390435802)
    ,@("N17", 1.103079095722038)
    ,@("O17", 1.832724773092338)
    ,@("B18", 1.258283553487615)
    ,@("C18", 0.3105798143529199)
    ,@("E18", 0.2750432858543661)
    ,@("F18", 1.554753667435236)
    ,@("G18", 0.3983272833906568)
    ,@("H18", 0.5454491894476661)
    ,@("I18", 0.5421361872083885)
    ,@("J18", 0.02828623653569196)
    ,@("L18", 0.6052386312672695)
    ,@("N18", 1.102963373909077)
    ,@("O18", 1.8325417498028)
    ,@("B19", 1.249736468587855)
    ,@("C19", 0.3100463231228332)
    ,@("E19", 0.2750438139205649)
    ,@("F19", 1.554114575104506)
    ,@("G19", 0.3981359518224679)
    ,@("H19", 0.545632104742765)
    ,@("I19", 0.5425346495918966)
    ,@("J19", 0.02815009834294102)
    ,@("L19", 0.603981200376893)
    ,@("N19", 1.102930837895897)
    ,@("O19", 1.832499881480572)
    ,@("B20", 1.288197004976439)
    ,@("C20", 0.312446284404615)
    ,@("E20", 0.2750481780442158)
    ,@("F20", 1.557046553899681)
    ,@("G20", 0.3990176380254979)
    ,@("H20", 0.5448264776581198)
    ,@("I20", 0.540763278167276)
    ,@("J20", 0.0287623366120755)
    ,@("L20", 0.6096532384291748)
    ,@("N20", 1.103103668126337)
    ,@("O20", 1.832768195287656)
    ,@("B21", 1.417370452026489)
    ,@("C21", 0.3204950210821664)
    ,@("E21", 0.2751805952116158)
    ,@("F21", 1.567876790830468)
    ,@("G21", 0.4023421655361261)
    ,@("H21", 0.5424260408903336)
    ,@("I21", 0.5351926521682486)
    ,@("J21", 0.03081228685948645)
    ,@("L21", 0.628944660417659)
    ,@("N21", 1.104145645475825)
    ,@("O21", 1.835068647077009)
    ,@("B22", 1.501746203376911)
    ,@("C22", 0.3257435121078629)
    ,@("E22", 0.2753546533768017)
    ,@("F22", 1.575683951949827)
    ,@("G22", 0.4047853635176324)
    ,@("H22", 0.5410854255245994)
    ,@("I22", 0.531836145315765)
    ,@("J22", 0.03214663470512846)
    ,@("L22", 0.6417255577415801)
    ,@("N22", 1.105169288582857)
    ,@("O22", 1.837615583892358)
    ,@("B23", 1.45671819684037)
    ,@("C23", 0.3229434261459261)
    ,@("E23", 0.275253763478954)
    ,@("F23", 1.57145047541681)
    ,@("G23", 0.4034566076168318)
    ,@("H23", 0.5417800641223067)
    ,@("I23", 0.5336015509252725)
    ,@("J23", 0.03143497190981748)
    ,@("L23", 0.6348884454443464)
    ,@("N23", 1.104591671962154)
    ,@("O23", 1.836160708491121)
    ,@("B24", 1.286085093026315)
    ,@("C24", 0.3123145436754413)
    ,@("E24", 0.2750474947603792)
    ,@("F24", 1.556881857791851)
    ,@("G24", 0.3989678587393115)
    ,@("H24", 0.5448695659919593)
    ,@("I24", 0.5408591210825016)
    ,@("J24", 0.02872874173694129)
    ,@("L24", 0.6093408720424236)
    ,@("N24", 1.103092438745691)
    ,@("O24", 1.832748199916949)
    ,@("B25", 1.10189573025599)
    ,@("C25", 0.3008038950728178)
    ,@("E25", 0.2752023392893399)
    ,@("F25", 1.544302237252978)
    ,@("G25", 0.3952838448637124)
    ,@("H25", 0.5491826875098553)
    ,@("I25", 0.5499059035159632)
    ,@("J25", 0.02578723447748743)
    ,@("L25", 0.5825364946146294)
    ,@("N25", 1.104131231087905)
    ,@("O25", 1.833541903917499)
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

Write-Host "Updated $($updates.Count) cells"